$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Populate the new "Examples" column (C2:C8) with the example sentences.
# ---------------------------------------------------------------------------
$examples = @(
    "I was just thinking about you when you called. What a coincidence!",
    "We haven't met in years. Let's sit down and catch up over coffee.",
    "I haven't seen you since graduation. What are you up to these days?",
    "You're late? No worries, the meeting hasn't started yet.",
    "I wouldn't ask about his salary; some people might find it a little personal.",
    "Do I want the red or blue shirt? Actually, I haven't made up my mind.",
    "Between the two candidates, I'm leaning toward Clancy because of his experience."
)
for ($i = 0; $i -lt $examples.Count; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $examples[$i]
}

# ---------------------------------------------------------------------------
# 2. Borders: thin black -> medium black around the whole A1:C8 table.
# ---------------------------------------------------------------------------
$allData = $ws.Range("A1:C8")
$allData.Borders.LineStyle = 1
$allData.Borders.Weight = -4138
$allData.Borders.Color = 0

# ---------------------------------------------------------------------------
# 3. Fonts: every data cell becomes 11pt Arial, color #1F1F1F.
#    Column A (the "label" column, which shares the header's bold look) is
#    bold; columns B & C are regular weight.
# ---------------------------------------------------------------------------
$headerRow = $ws.Range("A1:C1")
$headerRow.Font.Name = "Arial"
$headerRow.Font.Size = 11
$headerRow.Font.Color = 2039583
$headerRow.Font.Bold = $true

$colA = $ws.Range("A2:A8")
$colA.Font.Name = "Arial"
$colA.Font.Size = 11
$colA.Font.Color = 2039583
$colA.Font.Bold = $true

$colBC = $ws.Range("B2:C8")
$colBC.Font.Name = "Arial"
$colBC.Font.Size = 11
$colBC.Font.Color = 2039583

# ---------------------------------------------------------------------------
# 4. Alignment: left / center / wrap / indent 1 / reading order left-to-right
#    for every cell in the table.
# ---------------------------------------------------------------------------
$allData.HorizontalAlignment = -4131
$allData.VerticalAlignment = -4108
$allData.WrapText = $true
$allData.IndentLevel = 1
$allData.ReadingOrder = 1

# ---------------------------------------------------------------------------
# 5. Column widths (auto-fit result baked in explicitly).
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 41.7109375
$ws.Columns("C").ColumnWidth = 57.42578125

# ---------------------------------------------------------------------------
# 6. Row heights (auto-fit result baked in explicitly).
# ---------------------------------------------------------------------------
$ws.Rows(1).RowHeight = 30.75
$ws.Rows(2).RowHeight = 29.25
$ws.Rows(3).RowHeight = 29.25
$ws.Rows(4).RowHeight = 29.25
$ws.Rows(5).RowHeight = 15.75
$ws.Rows(6).RowHeight = 30.75
$ws.Rows(7).RowHeight = 29.25
$ws.Rows(8).RowHeight = 29.25

# ---------------------------------------------------------------------------
# 7. Selection moves to C13, mirroring the author's last selected cell.
# ---------------------------------------------------------------------------
$ws.Range("C13").Select()
